$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Countries list re-sorted: Belice / Nueva Caledonia swapped places (row 199 <-> row 200) ---
$ws.Range("A199").Value = "Nueva Caledonia"
$ws.Range("D199").Value = 18
$ws.Range("H199").Value = 0

$ws.Range("A200").Value = "Belice"
$ws.Range("D200").Value = 16
$ws.Range("H200").Value = 2

# --- Countries list re-sorted: Sahara Occidental / San Bartolome / Bonaire, San Eustaquio y Saba rotated (rows 214-216) ---
# Numeric columns are identical across these three rows, so only the country names move.
$ws.Range("A214").Value = "Bonaire, San Eustaquio y Saba"
$ws.Range("A215").Value = "Sahara Occidental"
$ws.Range("A216").Value = "San Bartolome"

# --- Update the "last refreshed" timestamp string (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 24 de Mayo de 2020 a las 10:35"

# --- Rumania (row 40) stats update ---
$ws.Range("E40").Value = 5491
$ws.Range("G40").Value = 3
$ws.Range("H40").Value = 1179

# --- Filipinas (row 46) stats update ---
$ws.Range("B46").Value = 14035
$ws.Range("C46").Value = 258
$ws.Range("D46").Value = 3249
$ws.Range("E46").Value = 9918
$ws.Range("G46").Value = 5
$ws.Range("H46").Value = 868
